$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has two header rows (row 1 and row 2) followed by
# the data rows (3..14). The edit collapses the two header rows into a
# single, new header row and adds two new leading "idx"/"idx2" columns'
# worth of headers plus renamed headers for the remaining columns.
# Deleting row 2 merges the two header rows into one (row 1) and shifts
# all the data rows up by one (so former row 3 becomes row 2, etc.),
# which matches the target layout (12 data rows, rows 2-13).
$ws.Rows.Item(2).Delete()

# Make sure the new header row cells start from a clean slate before
# writing the new header text, so columns A:E end up with the sheet's
# default (unstyled) formatting.
$ws.Range("A1:K1").ClearFormats()

$ws.Cells.Item(1,1).Value = "idx"
$ws.Cells.Item(1,2).Value = "idx2"
$ws.Cells.Item(1,3).Value = "Name"
$ws.Cells.Item(1,4).Value = "Date Start"
$ws.Cells.Item(1,5).Value = "Date End"

# Columns F:K keep using the same data font (Arial 9, like the rest of
# the table) while getting their new header captions.
$hdrRange = $ws.Range("F1:K1")
$hdrRange.Font.Name = "Arial"
$hdrRange.Font.Size = 9

$ws.Cells.Item(1,6).Value = "(m3/s)"
$ws.Cells.Item(1,7).Value = "(MW1)"
$ws.Cells.Item(1,8).Value = "(MW2)"
$ws.Cells.Item(1,9).Value = "(GWh) Winter"
$ws.Cells.Item(1,10).Value = "(GWh) Summer"
$ws.Cells.Item(1,11).Value = "(GWh) Year"

# Match the workbook's recorded selection after the edit.
$ws.Range("A2:K2").Select()
